$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F245").Value = 3501
$ws.Range("G245").Value = 79

$ws.Range("F272").Value = 30946
$ws.Range("G272").Value = 1671

$ws.Range("F273").Value = 31850
$ws.Range("G273").Value = 1674

$ws.Range("F278").Value = 30127
$ws.Range("G278").Value = 2076

$ws.Range("F279").Value = 42865
$ws.Range("G279").Value = 3059

$ws.Range("F280").Value = 34365
$ws.Range("G280").Value = 2298

$ws.Range("F281").Value = 45487
$ws.Range("G281").Value = 3151

$ws.Range("F282").Value = 47289
$ws.Range("G282").Value = 2830

$ws.Range("F285").Value = 41990
$ws.Range("G285").Value = 3423

$ws.Range("F286").Value = 55180
$ws.Range("G286").Value = 4288

$ws.Range("F287").Value = 59035
$ws.Range("G287").Value = 3733

$ws.Range("F288").Value = 58715
$ws.Range("G288").Value = 3932

$ws.Range("F289").Value = 63000
$ws.Range("G289").Value = 3654

$ws.Range("F292").Value = 82669
$ws.Range("G292").Value = 7305

$ws.Range("F293").Value = 82504
$ws.Range("G293").Value = 5769

$ws.Range("F294").Value = 93590
$ws.Range("G294").Value = 4946

$ws.Range("F295").Value = 17200

$ws.Range("F299").Value = 65431
$ws.Range("G299").Value = 6877

$ws.Range("F300").Value = 72179
$ws.Range("G300").Value = 6968

$ws.Range("F301").Value = 71799
$ws.Range("G301").Value = 5667

$ws.Range("F302").Value = 77713
$ws.Range("G302").Value = 5657

$ws.Range("F305").Value = 3268

$ws.Range("F306").Value = 74575
$ws.Range("G306").Value = 7598

$ws.Range("F307").Value = 76839
$ws.Range("G307").Value = 6358

$ws.Range("F308").Value = 15362
$ws.Range("G308").Value = 1050

$ws.Range("F309").Value = 77579
$ws.Range("G309").Value = 5532

$ws.Range("F310").Value = 78814
$ws.Range("G310").Value = 4058

$ws.Range("F311").Value = 61386
$ws.Range("G311").Value = 1928

$ws.Range("F313").Value = 75216
$ws.Range("G313").Value = 3452

$ws.Range("F314").Value = 63961
$ws.Range("G314").Value = 3146

$ws.Range("F315").Value = 56361
$ws.Range("G315").Value = 2656

$ws.Range("F316").Value = 50464
$ws.Range("G316").Value = 2299

$ws.Range("F317").Value = 63917
$ws.Range("G317").Value = 2186

$ws.Range("F318").Value = 48964
$ws.Range("G318").Value = 1135

$ws.Range("F319").Value = 41306

$ws.Range("F320").Value = 71150
$ws.Range("G320").Value = 3300

$ws.Range("F321").Value = 93114
$ws.Range("G321").Value = 2812

$ws.Range("F322").Value = 109648
$ws.Range("G322").Value = 2347

$ws.Range("F323").Value = 216156
$ws.Range("G323").Value = 3100

$ws.Range("F324").Value = 239826
$ws.Range("G324").Value = 2790

$ws.Range("F325").Value = 762456
$ws.Range("G325").Value = 6417

$ws.Range("F326").Value = 439336
$ws.Range("G326").Value = 3945

$ws.Range("F327").Value = 224301
$ws.Range("G327").Value = 2707

$ws.Range("F328").Value = 180524
$ws.Range("G328").Value = 2642

$ws.Range("F329").Value = 83001
$ws.Range("G329").Value = 1759

$ws.Range("F330").Value = 72428
$ws.Range("G330").Value = 2084

$ws.Range("F331").Value = 153111
$ws.Range("G331").Value = 2665

$ws.Range("F332").Value = 451849
$ws.Range("G332").Value = 4472

$ws.Range("F333").Value = 268613
$ws.Range("G333").Value = 2889
